# Auto-generated script to update Brynhildr_Profits market-data cells
# across all 8 worksheets, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 4619.2856
$ws.Range("I32").Value = 2196.25
$ws.Range("J32").Value = 6110.385
$ws.Range("K32").Value = 2196.25
$ws.Range("L32").Value = 6110.385
$ws.Range("M32").Value = -1870.25
$ws.Range("N32").Value = -6762.385
# Row 33
$ws.Range("H33").Value = 222.8
$ws.Range("I33").Value = 222.8
$ws.Range("K33").Value = 222.8
$ws.Range("M33").Value = 6.199999999999989
# Row 48
$ws.Range("H48").Value = 1551.6666
$ws.Range("I48").Value = 1200
$ws.Range("J48").Value = 1622
$ws.Range("K48").Value = 3600
$ws.Range("L48").Value = 4866
$ws.Range("M48").Value = -3308
$ws.Range("N48").Value = -5450
# Row 51
$ws.Range("H51").Value = 3131.1333
$ws.Range("I51").Value = 2909.7805
$ws.Range("J51").Value = 5400
$ws.Range("K51").Value = 2909.7805
$ws.Range("L51").Value = 5400
$ws.Range("M51").Value = -2425.7805
$ws.Range("N51").Value = -6368
# Row 53
$ws.Range("H53").Value = 288.4737
$ws.Range("I53").Value = 304.16666
$ws.Range("J53").Value = 261.57144
$ws.Range("K53").Value = 304.16666
$ws.Range("L53").Value = 261.57144
$ws.Range("M53").Value = 332.83334
$ws.Range("N53").Value = -1535.57144
# Row 56
$ws.Range("H56").Value = 1551.6666
$ws.Range("I56").Value = 1200
$ws.Range("J56").Value = 1622
$ws.Range("K56").Value = 3600
$ws.Range("L56").Value = 4866
$ws.Range("M56").Value = -3066
$ws.Range("N56").Value = -5934
# Row 64
$ws.Range("H64").Value = 3571.4285
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496
# Row 67
$ws.Range("H67").Value = 3571.4285
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716
# Row 98
$ws.Range("H98").Value = 1291.1613
$ws.Range("I98").Value = 1352.6207
$ws.Range("K98").Value = 1352.6207
$ws.Range("M98").Value = 145.3793000000001
# Row 107
$ws.Range("H107").Value = 3431.8635
$ws.Range("I107").Value = 3431.8635
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3431.8635
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1511.8635
$ws.Range("N107").ClearContents()
# Row 122
$ws.Range("H122").Value = 1291.1613
$ws.Range("I122").Value = 1352.6207
$ws.Range("K122").Value = 4057.8621
$ws.Range("M122").Value = -1607.8621
# Row 132
$ws.Range("H132").Value = 1518.3334
$ws.Range("I132").Value = 1518.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4555.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2025.0002
$ws.Range("N132").ClearContents()
# Row 135
$ws.Range("H135").Value = 1795.2858
$ws.Range("I135").Value = 306.2
$ws.Range("K135").Value = 2755.8
$ws.Range("M135").Value = -220.7999999999997
# Row 141
$ws.Range("H141").Value = 1349.5
$ws.Range("I141").Value = 1349.5
$ws.Range("K141").Value = 4048.5
$ws.Range("M141").Value = 1131.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 21331
$ws.Range("I2").Value = 18497
$ws.Range("K2").Value = 18497
$ws.Range("M2").Value = -18384
# Row 45
$ws.Range("H45").Value = 4648.3
$ws.Range("I45").Value = 3721
$ws.Range("J45").Value = 5266.5
$ws.Range("K45").Value = 3721
$ws.Range("L45").Value = 5266.5
$ws.Range("M45").Value = -3344
$ws.Range("N45").Value = -6020.5
# Row 74
$ws.Range("H74").Value = 8512.527
$ws.Range("I74").Value = 3779
$ws.Range("J74").Value = 15139.467
$ws.Range("K74").Value = 3779
$ws.Range("L74").Value = 15139.467
$ws.Range("M74").Value = -2905
$ws.Range("N74").Value = -16887.467
# Row 77
$ws.Range("H77").Value = 8512.527
$ws.Range("I77").Value = 3779
$ws.Range("J77").Value = 15139.467
$ws.Range("K77").Value = 18895
$ws.Range("L77").Value = 75697.33500000001
$ws.Range("M77").Value = -14527
$ws.Range("N77").Value = -84433.33500000001
# Row 116
$ws.Range("H116").Value = 21331
$ws.Range("I116").Value = 18497
$ws.Range("K116").Value = 18497
$ws.Range("M116").Value = -16203
# Row 122
$ws.Range("H122").Value = 2432.9167
$ws.Range("I122").Value = 2219.5
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 6658.5
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -4208.5
$ws.Range("N122").Value = -15400
# Row 132
$ws.Range("H132").Value = 4119.304
$ws.Range("I132").Value = 3455.647
$ws.Range("K132").Value = 10366.941
$ws.Range("M132").Value = -7836.940999999999
# Row 133
$ws.Range("H133").Value = 78813
$ws.Range("J133").Value = 78813
$ws.Range("L133").Value = 78813
$ws.Range("N133").Value = -83873

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 21331
$ws.Range("I3").Value = 18497
$ws.Range("K3").Value = 18497
$ws.Range("M3").Value = -18383
# Row 55
$ws.Range("H55").Value = 65000
$ws.Range("J55").Value = 65000
$ws.Range("L55").Value = 65000
$ws.Range("N55").Value = -65546
# Row 107
$ws.Range("H107").Value = 1849.2858
$ws.Range("I107").Value = 1889
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 1889
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 31
$ws.Range("N107").Value = -5590
# Row 134
$ws.Range("H134").Value = 6181.974
$ws.Range("I134").Value = 3332.5625
$ws.Range("K134").Value = 9997.6875
$ws.Range("M134").Value = -7462.6875
# Row 141
$ws.Range("H141").Value = 349999.5
$ws.Range("J141").Value = 349999.5
$ws.Range("L141").Value = 349999.5
$ws.Range("N141").Value = -360359.5

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 5839.4043
$ws.Range("I58").Value = 4590.795
$ws.Range("J58").Value = 11926.375
$ws.Range("K58").Value = 4590.795
$ws.Range("L58").Value = 11926.375
$ws.Range("M58").Value = -4387.795
$ws.Range("N58").Value = -12332.375
# Row 70
$ws.Range("H70").Value = 40744.5
$ws.Range("J70").Value = 40744.5
$ws.Range("L70").Value = 40744.5
$ws.Range("N70").Value = -41374.5
# Row 73
$ws.Range("H73").Value = 40744.5
$ws.Range("J73").Value = 40744.5
$ws.Range("L73").Value = 40744.5
$ws.Range("N73").Value = -42928.5
# Row 107
$ws.Range("H107").Value = 957.75
$ws.Range("I107").Value = 753.2222
$ws.Range("J107").Value = 1220.7142
$ws.Range("K107").Value = 753.2222
$ws.Range("L107").Value = 1220.7142
$ws.Range("M107").Value = 1166.7778
$ws.Range("N107").Value = -5060.7142
# Row 122
$ws.Range("H122").Value = 22516.68
$ws.Range("I122").Value = 2805.6667
$ws.Range("K122").Value = 8417.000100000001
$ws.Range("M122").Value = -5967.000100000001
# Row 132
$ws.Range("H132").Value = 4215.0884
$ws.Range("I132").Value = 4023.7666
$ws.Range("J132").Value = 5650
$ws.Range("K132").Value = 12071.2998
$ws.Range("L132").Value = 16950
$ws.Range("M132").Value = -9541.299800000001
$ws.Range("N132").Value = -22010
# Row 134
$ws.Range("H134").Value = 1418.644
$ws.Range("I134").Value = 1418.644
$ws.Range("K134").Value = 4255.932
$ws.Range("M134").Value = -1720.932
# Row 136
$ws.Range("H136").Value = 5839.4043
$ws.Range("I136").Value = 4590.795
$ws.Range("J136").Value = 11926.375
$ws.Range("K136").Value = 13772.385
$ws.Range("L136").Value = 35779.125
$ws.Range("M136").Value = -11222.385
$ws.Range("N136").Value = -40879.125

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 8400338
$ws.Range("I4").Value = 8400338
$ws.Range("K4").Value = 25201014
$ws.Range("M4").Value = -25200902
# Row 15
$ws.Range("H15").Value = 33.75
$ws.Range("I15").Value = 35.75
$ws.Range("K15").Value = 107.25
$ws.Range("M15").Value = 32.75
# Row 122
$ws.Range("H122").Value = 1614795.8
$ws.Range("I122").Value = 5376403
$ws.Range("K122").Value = 48387627
$ws.Range("M122").Value = -48385177
# Row 131
$ws.Range("H131").Value = 1936.4584
$ws.Range("J131").Value = 2021.591
$ws.Range("L131").Value = 6064.772999999999
$ws.Range("N131").Value = -16144.773
# Row 132
$ws.Range("H132").Value = 3232.7778
$ws.Range("J132").Value = 3033.3333
$ws.Range("L132").Value = 27299.9997
$ws.Range("N132").Value = -32359.9997

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1996.4445
$ws.Range("J113").Value = 1987.5
$ws.Range("L113").Value = 1987.5
$ws.Range("N113").Value = -6327.5
# Row 132
$ws.Range("H132").Value = 10171.412
$ws.Range("I132").Value = 10859.333
$ws.Range("K132").Value = 32577.999
$ws.Range("M132").Value = -30047.999

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 4290.7827
$ws.Range("J46").Value = 4666
$ws.Range("L46").Value = 4666
$ws.Range("N46").Value = -5042
# Row 55
$ws.Range("H55").Value = 1388.3715
$ws.Range("I55").Value = 1414.6923
$ws.Range("J55").Value = 1372.8182
$ws.Range("K55").Value = 1414.6923
$ws.Range("L55").Value = 1372.8182
$ws.Range("M55").Value = -1241.6923
$ws.Range("N55").Value = -1718.8182
# Row 74
$ws.Range("H74").Value = 46296.89
$ws.Range("J74").Value = 52259
$ws.Range("L74").Value = 52259
$ws.Range("N74").Value = -54255
# Row 77
$ws.Range("H77").Value = 46296.89
$ws.Range("J77").Value = 52259
$ws.Range("L77").Value = 156777
$ws.Range("N77").Value = -166761
# Row 122
$ws.Range("H122").Value = 4800.9287
$ws.Range("I122").Value = 3571
$ws.Range("J122").Value = 6440.8335
$ws.Range("K122").Value = 10713
$ws.Range("L122").Value = 19322.5005
$ws.Range("M122").Value = -8263
$ws.Range("N122").Value = -24222.5005
# Row 132
$ws.Range("H132").Value = 2777.875
$ws.Range("I132").Value = 2777.875
$ws.Range("K132").Value = 8333.625
$ws.Range("M132").Value = -5803.625

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 4260.3335
$ws.Range("I132").Value = 4023.8333
$ws.Range("K132").Value = 12071.4999
$ws.Range("M132").Value = -9541.499899999999
